$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.939.13"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").Value = "3.764.15"
$ws.Range("E3").Value = "  -1.82%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'629.45"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("D6").Value = "'165.63"
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").Value = "3.764.28"
$ws.Range("E7").Value = "  -1.75%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("E10").Value = "  -2.36%  "
$ws.Range("E11").Value = "  +0.36%  "
$ws.Range("E12").Value = "  +0.63%  "
$ws.Range("E13").Value = "  -4.57%  "
$ws.Range("D14").Value = "'34.84"
$ws.Range("E14").Value = "  -3.07%  "
$ws.Range("D15").Value = "4.399.32"
$ws.Range("E15").Value = "  -1.72%  "
$ws.Range("D16").Value = "3.757.67"
$ws.Range("E16").Value = "  -2.03%  "
$ws.Range("D17").Value = "69.000.88"
$ws.Range("E17").Value = "  -0.22%  "
$ws.Range("E18").Value = "  -3.39%  "
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("D20").Value = "'7.01"
$ws.Range("E20").Value = "  -2.06%  "
$ws.Range("D21").Value = "'461.72"
$ws.Range("E21").Value = "  -1.23%  "
$ws.Range("D22").Value = "'9.48"
$ws.Range("E22").Value = "  -3.04%  "
$ws.Range("D23").Value = "'0.703"
$ws.Range("E23").Value = "  -0.86%  "
$ws.Range("E24").Value = "  -6.06%  "
$ws.Range("D25").Value = "'82.03"
$ws.Range("E25").Value = "  -2.43%  "
$ws.Range("D26").Value = "'12.11"
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("E27").Value = "  -1.94%  "
$ws.Range("E28").Value = "  +0.38%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").Value = "3.921.24"
$ws.Range("E30").Value = "  -1.52%  "
$ws.Range("E31").Value = "  +1.03%  "
$ws.Range("E32").Value = "  -0.19%  "
$ws.Range("D33").Value = "'7.05"
$ws.Range("E33").Value = "  -3.79%  "
$ws.Range("D34").Value = "'28.30"
$ws.Range("E34").Value = "  -3.27%  "
$ws.Range("E35").Value = "  +17.35%  "
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  +0.16%  "
$ws.Range("D37").Value = "3.718.09"
$ws.Range("E37").Value = "  -1.62%  "
$ws.Range("E38").Value = "  -2.65%  "
$ws.Range("E39").Value = "  -1.65%  "
$ws.Range("E40").Value = "  +0.53%  "
$ws.Range("D41").Value = "'5.78"
$ws.Range("E41").Value = "  -2.48%  "
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  +0.15%  "
$ws.Range("D43").Value = "'0.959"
$ws.Range("E43").Value = "  -2.16%  "
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").Value = "'156.77"
$ws.Range("E45").Value = "  -0.28%  "
$ws.Range("E46").Value = "  +3.62%  "
$ws.Range("D47").Value = "'1.41"
$ws.Range("E47").Value = "  -0.43%  "
$ws.Range("D48").Value = "'46.97"
$ws.Range("E48").Value = "  +0.09%  "
$ws.Range("D49").Value = "'42.75"
$ws.Range("E49").Value = "  -0.68%  "
$ws.Range("E50").Value = "  -2.65%  "
$ws.Range("D51").Value = "'8.33"
$ws.Range("E51").Value = "  -1.18%  "
